$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new empty column before column A.
# This shifts: A->B, B->C, C->D, D->E, E->F, F->G
$ws.Columns.Item(1).Insert()

# --- Header row ---
# B1 (new) gets the "segments" header, styled like the other header cells.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B1").Value = "segments"

# --- Data rows (2-20) ---
# Column A gets numeric segment indices (0-based), styled like the old
# label column (bold, centered, bordered).
# Column B keeps the label text (shifted from former column A) but loses
# that bold/bordered style, taking on the plain default style instead
# (copied from the neighboring value cell in column C).
$segments = @(
    "background",
    "back_bumper",
    "back_glass",
    "back_left_door",
    "back_left_light",
    "back_right_door",
    "back_right_light",
    "front_bumper",
    "front_glass",
    "front_left_door",
    "front_left_light",
    "front_right_door",
    "front_right_light",
    "hood",
    "left_mirror",
    "right_mirror",
    "tailgate",
    "trunk",
    "wheel"
)

for ($i = 0; $i -lt $segments.Length; $i++) {
    $row = $i + 2

    # Style + value for the new index column A
    $ws.Range("B$row").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range("A$row").Value = $i

    # Restyle the shifted label column B to the plain/default look
    $ws.Range("C$row").Copy()
    $ws.Range("B$row").PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range("B$row").Value = $segments[$i]
}

$excel.CutCopyMode = $false
